# "Many minor bug fixes"
#
# 1. DATASET_OPTIONS sheet: the 'batch' column (col G) is removed entirely.
# 2. DATASET_OPTIONS sheet: the 'enhance_signals' values (col C) for the
#    3-seq and ChIP-seq rows were stored as a native boolean TRUE; they are
#    changed to the plain text string "TRUE" (matching the other rows which
#    already store FALSE/TRUE as text).
# 3. GLOBAL_OPTIONS sheet: track_height_cm value changes from 0.4 to 0.3.
# 4. The SAMPLES tab becomes the active/selected tab (instead of
#    DATASET_OPTIONS).

$wb = $excel.ActiveWorkbook

$wsDataset = $wb.Worksheets.Item("DATASET_OPTIONS")

# --- Remove the 'batch' column (column G) ---
$wsDataset.Columns.Item(7).Delete()

# --- Fix the enhance_signals cells that were stored as native booleans ---
# Using a leading apostrophe forces Excel to store the value as literal text
# (shared string) rather than re-interpreting "TRUE" as a boolean again.
$wsDataset.Cells.Item(2,3).Value = "'TRUE"
$wsDataset.Cells.Item(2,3).Style = "Normal"
$wsDataset.Cells.Item(3,3).Value = "'TRUE"
$wsDataset.Cells.Item(3,3).Style = "Normal"

# --- GLOBAL_OPTIONS: track_height_cm 0.4 -> 0.3 ---
$wsGlobal = $wb.Worksheets.Item("GLOBAL_OPTIONS")
$wsGlobal.Cells.Item(69,2).Value = "'0.3"
$wsGlobal.Cells.Item(69,2).Style = "Normal"

# --- Make SAMPLES the active sheet/tab ---
$wsSamples = $wb.Worksheets.Item("SAMPLES")
$wsSamples.Activate()
$wsSamples.Select()

$wb.Save()
